$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 to make room for the new earliest data point (2007),
# shifting the existing data rows down by one.
$ws.Rows("2:2").Insert()

# The inserted row copies formatting from the row above (the header row) across
# all five columns; reset that and re-apply the date style only to column A
# (matching the rest of the table: bold, centered, bordered, custom date format).
$ws.Range("A2:E2").ClearFormats()
$ws.Range("A2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Borders.LineStyle = 1

# Write out the full corrected data table (dates, years, and recomputed
# y_0_forecast / y_1_forecast values) for rows 2:19.
$ws.Range("A2").Value2 = 39400
$ws.Range("B2").Value2 = 2007
$ws.Range("C2").Value2 = 0.4235526809466261
$ws.Range("D2").Value2 = 2008
$ws.Range("E2").Value2 = 1.652713177926435
$ws.Range("A3").Value2 = 39765
$ws.Range("B3").Value2 = 2008
$ws.Range("C3").Value2 = -0.5718076928962645
$ws.Range("D3").Value2 = 2009
$ws.Range("E3").Value2 = 0.4618648366506939
$ws.Range("A4").Value2 = 40130
$ws.Range("B4").Value2 = 2009
$ws.Range("C4").Value2 = 0.3486139762225005
$ws.Range("D4").Value2 = 2010
$ws.Range("E4").Value2 = -2.340608900318997
$ws.Range("A5").Value2 = 40494
$ws.Range("B5").Value2 = 2010
$ws.Range("C5").Value2 = -0.1384957661262898
$ws.Range("D5").Value2 = 2011
$ws.Range("E5").Value2 = 1.659950937631938
$ws.Range("A6").Value2 = 40862
$ws.Range("B6").Value2 = 2011
$ws.Range("C6").Value2 = 1.566479473280147
$ws.Range("D6").Value2 = 2012
$ws.Range("E6").Value2 = 2.497560574935442
$ws.Range("A7").Value2 = 41228
$ws.Range("B7").Value2 = 2012
$ws.Range("C7").Value2 = 0.7307568962936939
$ws.Range("D7").Value2 = 2013
$ws.Range("E7").Value2 = 1.029202372425875
$ws.Range("A8").Value2 = 41592
$ws.Range("B8").Value2 = 2013
$ws.Range("C8").Value2 = 0.818818812164257
$ws.Range("D8").Value2 = 2014
$ws.Range("E8").Value2 = 0.7004752402133052
$ws.Range("A9").Value2 = 41957
$ws.Range("B9").Value2 = 2014
$ws.Range("C9").Value2 = 0.9180054319587239
$ws.Range("D9").Value2 = 2015
$ws.Range("E9").Value2 = 2.577405783391451
$ws.Range("A10").Value2 = 42321
$ws.Range("B10").Value2 = 2015
$ws.Range("C10").Value2 = 1.984684278296656
$ws.Range("D10").Value2 = 2016
$ws.Range("E10").Value2 = 2.047428048848832
$ws.Range("A11").Value2 = 42689
$ws.Range("B11").Value2 = 2016
$ws.Range("C11").Value2 = 1.755995812646982
$ws.Range("D11").Value2 = 2017
$ws.Range("E11").Value2 = 1.55296524673576
$ws.Range("A12").Value2 = 43053
$ws.Range("B12").Value2 = 2017
$ws.Range("C12").Value2 = 1.946965557828384
$ws.Range("D12").Value2 = 2018
$ws.Range("E12").Value2 = 0.2326081529569146
$ws.Range("A13").Value2 = 43418
$ws.Range("B13").Value2 = 2018
$ws.Range("C13").Value2 = 1.06432145354225
$ws.Range("D13").Value2 = 2019
$ws.Range("E13").Value2 = -0.6993904531539141
$ws.Range("A14").Value2 = 43783
$ws.Range("B14").Value2 = 2019
$ws.Range("C14").Value2 = 1.361817904277696
$ws.Range("D14").Value2 = 2020
$ws.Range("E14").Value2 = 1.656460003703519
$ws.Range("A15").Value2 = 44159
$ws.Range("B15").Value2 = 2020
$ws.Range("C15").Value2 = -4.352425014431304
$ws.Range("D15").Value2 = 2021
$ws.Range("E15").Value2 = -0.9164352768978534
$ws.Range("A16").Value2 = 44525
$ws.Range("B16").Value2 = 2021
$ws.Range("C16").Value2 = -1.761645650979182
$ws.Range("D16").Value2 = 2022
$ws.Range("E16").Value2 = 2.211960525313206
$ws.Range("A17").Value2 = 44890
$ws.Range("B17").Value2 = 2022
$ws.Range("C17").Value2 = 5.20787683103745
$ws.Range("D17").Value2 = 2023
$ws.Range("E17").Value2 = 0.1996286546173343
$ws.Range("A18").Value2 = 45254
$ws.Range("B18").Value2 = 2023
$ws.Range("C18").Value2 = -0.9008525709169546
$ws.Range("D18").Value2 = 2024
$ws.Range("E18").Value2 = -0.01592365107300209
$ws.Range("A19").Value2 = 45618
$ws.Range("B19").Value2 = 2024
$ws.Range("C19").Value2 = 0.2738544794132824
$ws.Range("D19").Value2 = 2025
$ws.Range("E19").Value2 = 0.7861875694807674
